$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the measured values in column C (rows 2-4) ---
$ws.Range("C2").Value = 11
$ws.Range("C3").Value = 9.5
$ws.Range("C4").Value = 1.4

# --- Resize columns A and C to fixed widths (no longer auto "best fit") ---
# Excel's ColumnWidth is in "characters" and gets stored in the file as
# characters + 5/MaxDigitWidth; compensate so the saved <col width="..">
# lands on the desired value (27 and 27.25 respectively).
$ws.Columns.Item(1).ColumnWidth = 27 - 5 / 7
$ws.Columns.Item(3).ColumnWidth = 27.25 - 5 / 7

# --- Move the active selection to B3 ---
$ws.Range("B3").Select()
